# Applies scheduled-runner value updates to the Phoenix_Profits market-board sheets.
# Each block targets one worksheet/row; only currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are touched, matching the upstream diff.
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 927.5294
$ws.Range("I80").Value = 573.75
$ws.Range("J80").Value = 1242
$ws.Range("K80").Value = 1721.25
$ws.Range("L80").Value = 3726
$ws.Range("M80").Value = -723.25
$ws.Range("N80").Value = -5722

# Sheet ALC, row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 927.5294
$ws.Range("I83").Value = 573.75
$ws.Range("J83").Value = 1242
$ws.Range("K83").Value = 5163.75
$ws.Range("L83").Value = 11178
$ws.Range("M83").Value = -171.75
$ws.Range("N83").Value = -21162

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3000
$ws.Range("I111").Value = 3000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 9000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -5933

# Sheet ALC, row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2345.4443
$ws.Range("I127").Value = 2345.4443
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 7036.3329
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -2076.3329

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1152432.5
$ws.Range("I132").Value = 1152432.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3457297.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3454767.5

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1479.1111
$ws.Range("I137").Value = 1405.5333
$ws.Range("J137").Value = 1847
$ws.Range("K137").Value = 4216.5999
$ws.Range("L137").Value = 5541
$ws.Range("M137").Value = -1666.5999
$ws.Range("N137").Value = -10641

# Sheet ALC, row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1812.2632
$ws.Range("I2").Value = 1856.2667
$ws.Range("J2").Value = 1647.25
$ws.Range("K2").Value = 1856.2667
$ws.Range("L2").Value = 1647.25
$ws.Range("M2").Value = -1743.2667
$ws.Range("N2").Value = -1873.25

# Sheet ARM, row 64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 58999.555
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 58999.555
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 58999.555
$ws.Range("N64").Value = -59495.555

# Sheet ARM, row 67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 58999.555
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 58999.555
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 58999.555
$ws.Range("N67").Value = -60715.555

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3099.12
$ws.Range("I74").Value = 2083.7856
$ws.Range("J74").Value = 4391.364
$ws.Range("K74").Value = 2083.7856
$ws.Range("L74").Value = 4391.364
$ws.Range("M74").Value = -1209.7856
$ws.Range("N74").Value = -6139.364

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3099.12
$ws.Range("I77").Value = 2083.7856
$ws.Range("J77").Value = 4391.364
$ws.Range("K77").Value = 10418.928
$ws.Range("L77").Value = 21956.82
$ws.Range("M77").Value = -6050.928
$ws.Range("N77").Value = -30692.82

# Sheet ARM, row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1659.5454
$ws.Range("I88").Value = 1785.6
$ws.Range("J88").Value = 1554.5
$ws.Range("K88").Value = 1785.6
$ws.Range("L88").Value = 1554.5
$ws.Range("M88").Value = -1379.6
$ws.Range("N88").Value = -2366.5

# Sheet ARM, row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1659.5454
$ws.Range("I91").Value = 1785.6
$ws.Range("J91").Value = 1554.5
$ws.Range("K91").Value = 1785.6
$ws.Range("L91").Value = 1554.5
$ws.Range("M91").Value = -381.5999999999999
$ws.Range("N91").Value = -4362.5

# Sheet ARM, row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 22569
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 22569
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 22569
$ws.Range("N104").Value = -29557

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1812.2632
$ws.Range("I116").Value = 1856.2667
$ws.Range("J116").Value = 1647.25
$ws.Range("K116").Value = 1856.2667
$ws.Range("L116").Value = 1647.25
$ws.Range("M116").Value = 437.7333000000001
$ws.Range("N116").Value = -6235.25

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1812.2632
$ws.Range("I3").Value = 1856.2667
$ws.Range("J3").Value = 1647.25
$ws.Range("K3").Value = 1856.2667
$ws.Range("L3").Value = 1647.25
$ws.Range("M3").Value = -1742.2667
$ws.Range("N3").Value = -1875.25

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 43479228
$ws.Range("I86").Value = 76923990
$ws.Range("J86").Value = 1030.1
$ws.Range("K86").Value = 76923990
$ws.Range("L86").Value = 1030.1
$ws.Range("M86").Value = -76922867
$ws.Range("N86").Value = -3276.1

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 43479228
$ws.Range("I89").Value = 76923990
$ws.Range("J89").Value = 1030.1
$ws.Range("K89").Value = 384619950
$ws.Range("L89").Value = 5150.5
$ws.Range("M89").Value = -384614334
$ws.Range("N89").Value = -16382.5

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5102972.5
$ws.Range("I94").Value = 6410999.5
$ws.Range("J94").Value = 1667.5
$ws.Range("K94").Value = 6410999.5
$ws.Range("L94").Value = 1667.5
$ws.Range("M94").Value = -6410548.5
$ws.Range("N94").Value = -2569.5

# Sheet BSM, row 106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 37868.2
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 37868.2
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 37868.2
$ws.Range("N106").Value = -40392.2

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 40199
$ws.Range("I107").Value = 52395.934
$ws.Range("J107").Value = 3608.2
$ws.Range("K107").Value = 52395.934
$ws.Range("L107").Value = 3608.2
$ws.Range("M107").Value = -50475.934
$ws.Range("N107").Value = -7448.2

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1664.6
$ws.Range("I31").Value = 1115.4615
$ws.Range("J31").Value = 2259.5
$ws.Range("K31").Value = 1115.4615
$ws.Range("L31").Value = 2259.5
$ws.Range("M31").Value = -820.4614999999999
$ws.Range("N31").Value = -2849.5

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1664.6
$ws.Range("I34").Value = 1115.4615
$ws.Range("J34").Value = 2259.5
$ws.Range("K34").Value = 1115.4615
$ws.Range("L34").Value = 2259.5
$ws.Range("M34").Value = -913.4614999999999
$ws.Range("N34").Value = -2663.5

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 53849644
$ws.Range("I58").Value = 33336580
$ws.Range("J58").Value = 71432264
$ws.Range("K58").Value = 33336580
$ws.Range("L58").Value = 71432264
$ws.Range("M58").Value = -33336377
$ws.Range("N58").Value = -71432670

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3595.9443
$ws.Range("I99").Value = 3721.2
$ws.Range("J99").Value = 2969.6667
$ws.Range("K99").Value = 3721.2
$ws.Range("L99").Value = 2969.6667
$ws.Range("M99").Value = -2223.2
$ws.Range("N99").Value = -5965.6667

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3595.9443
$ws.Range("I126").Value = 3721.2
$ws.Range("J126").Value = 2969.6667
$ws.Range("K126").Value = 11163.6
$ws.Range("L126").Value = 8909.000100000001
$ws.Range("M126").Value = -8693.599999999999
$ws.Range("N126").Value = -13849.0001

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 53849644
$ws.Range("I136").Value = 33336580
$ws.Range("J136").Value = 71432264
$ws.Range("K136").Value = 100009740
$ws.Range("L136").Value = 214296792
$ws.Range("M136").Value = -100007190
$ws.Range("N136").Value = -214301892

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 76925300
$ws.Range("I113").Value = 350
$ws.Range("J113").Value = 83335704
$ws.Range("K113").Value = 1050
$ws.Range("L113").Value = 250007112
$ws.Range("M113").Value = 1120
$ws.Range("N113").Value = -250011452

# Sheet CUL, row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 18809.666
$ws.Range("I116").Value = 48029
$ws.Range("J116").Value = 4200
$ws.Range("K116").Value = 144087
$ws.Range("L116").Value = 12600
$ws.Range("M116").Value = -140645
$ws.Range("N116").Value = -19484

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 16682
$ws.Range("I113").Value = 20576
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 20576
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -18406
$ws.Range("N113").Value = -9340

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 47622828
$ws.Range("I126").Value = 83336270
$ws.Range("J126").Value = 4901.778
$ws.Range("K126").Value = 250008810
$ws.Range("L126").Value = 14705.334
$ws.Range("M126").Value = -250006340
$ws.Range("N126").Value = -19645.334

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2356743.2
$ws.Range("I132").Value = 3080358.5
$ws.Range("J132").Value = 4993.75
$ws.Range("K132").Value = 9241075.5
$ws.Range("L132").Value = 14981.25
$ws.Range("M132").Value = -9238545.5
$ws.Range("N132").Value = -20041.25

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4261.579
$ws.Range("I40").Value = 4122.875
$ws.Range("J40").Value = 5001.3335
$ws.Range("K40").Value = 4122.875
$ws.Range("L40").Value = 5001.3335
$ws.Range("M40").Value = -3986.875
$ws.Range("N40").Value = -5273.3335

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10274.762
$ws.Range("I132").Value = 9185.8125
$ws.Range("J132").Value = 13759.4
$ws.Range("K132").Value = 27557.4375
$ws.Range("L132").Value = 41278.2
$ws.Range("M132").Value = -25027.4375
$ws.Range("N132").Value = -46338.2

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2900.0833
$ws.Range("I132").Value = 2971.1428
$ws.Range("J132").Value = 2800.6
$ws.Range("K132").Value = 8913.4284
$ws.Range("L132").Value = 8401.8
$ws.Range("M132").Value = -6383.428400000001
$ws.Range("N132").Value = -13461.8
